$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65: new entry
$ws.Cells.Item(65, 1).Value = 202230988
$ws.Cells.Item(65, 2).Value = "https://vk.com/ekimovalex"
$ws.Cells.Item(65, 3).Value = "http://psytests.org/result?v=pfqOCrZZApR3RIR6CYBZ"
$ws.Cells.Item(65, 4).Value = 30
$ws.Cells.Item(65, 5).Value = "м"
$ws.Cells.Item(65, 7).Value = "ni"

# Row 66: new entry
$ws.Cells.Item(66, 1).Value = 175676662
$ws.Cells.Item(66, 2).Value = "https://vk.com/id175676662"
$ws.Cells.Item(66, 3).Value = "http://psytests.org/result?v=pfqKFHO9b-WOuhiwaByt"
$ws.Cells.Item(66, 4).Value = 24
$ws.Cells.Item(66, 5).Value = "ж"
$ws.Cells.Item(66, 7).Value = "ni"

# Update selection / view state to match the final selection on G66
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("G66").Select()
